$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add rows 429-435 of Landscaping Data (weather/growth observations for 2025-07-10)
$ws.Range("A429").Value = 45848
$ws.Range("A429").NumberFormat = "m/d/yy"
$ws.Range("B429").Value = "Flowering"
$ws.Range("C429").Value = "Large"
$ws.Range("D429").Value = 69
$ws.Range("E429").Value = 84
$ws.Range("G429").Value = 0
$ws.Range("H429").Value = 0.1
$ws.Range("I429").Value = "No"
$ws.Range("J429").Value = 2
$ws.Range("K429").Value = "Neutral"
$ws.Range("L429").Value = 7
$ws.Range("M429").Value = 0.63
$ws.Range("N429").Value = 68
$ws.Range("O429").Value = 30.03
$ws.Range("P429").Value = 4
$ws.Range("Q429").Value = 0.31
$ws.Range("R429").Value = 9.9
$ws.Range("S429").Value = 57
$ws.Range("T429").Value = 0
$ws.Range("A430").Value = 45848
$ws.Range("A430").NumberFormat = "m/d/yy"
$ws.Range("B430").Value = "Nonflowering"
$ws.Range("C430").Value = "Medium"
$ws.Range("D430").Value = 69
$ws.Range("E430").Value = 84
$ws.Range("G430").Value = 0
$ws.Range("H430").Value = 0.1
$ws.Range("I430").Value = "No"
$ws.Range("J430").Value = 3
$ws.Range("K430").Value = "Bright"
$ws.Range("L430").Value = 7
$ws.Range("M430").Value = 0.63
$ws.Range("N430").Value = 68
$ws.Range("O430").Value = 30.03
$ws.Range("P430").Value = 4
$ws.Range("Q430").Value = 0.31
$ws.Range("R430").Value = 9.9
$ws.Range("S430").Value = 57
$ws.Range("T430").Value = 0
$ws.Range("A431").Value = 45848
$ws.Range("A431").NumberFormat = "m/d/yy"
$ws.Range("B431").Value = "Nonflowering"
$ws.Range("C431").Value = "Small"
$ws.Range("D431").Value = 69
$ws.Range("E431").Value = 84
$ws.Range("G431").Value = 0
$ws.Range("H431").Value = 0
$ws.Range("I431").Value = "No"
$ws.Range("J431").Value = 3
$ws.Range("K431").Value = "Bright"
$ws.Range("L431").Value = 7
$ws.Range("M431").Value = 0.63
$ws.Range("N431").Value = 68
$ws.Range("O431").Value = 30.03
$ws.Range("P431").Value = 4
$ws.Range("Q431").Value = 0.31
$ws.Range("R431").Value = 9.9
$ws.Range("S431").Value = 57
$ws.Range("T431").Value = 0
$ws.Range("A432").Value = 45848
$ws.Range("A432").NumberFormat = "m/d/yy"
$ws.Range("B432").Value = "Nonflowering"
$ws.Range("C432").Value = "Medium"
$ws.Range("D432").Value = 69
$ws.Range("E432").Value = 84
$ws.Range("G432").Value = 0
$ws.Range("H432").Value = 0
$ws.Range("I432").Value = "No"
$ws.Range("J432").Value = 3
$ws.Range("K432").Value = "Neutral"
$ws.Range("L432").Value = 7
$ws.Range("M432").Value = 0.63
$ws.Range("N432").Value = 68
$ws.Range("O432").Value = 30.03
$ws.Range("P432").Value = 4
$ws.Range("Q432").Value = 0.31
$ws.Range("R432").Value = 9.9
$ws.Range("S432").Value = 57
$ws.Range("T432").Value = 0
$ws.Range("A433").Value = 45848
$ws.Range("A433").NumberFormat = "m/d/yy"
$ws.Range("B433").Value = "Nonflowering"
$ws.Range("C433").Value = "Medium"
$ws.Range("D433").Value = 69
$ws.Range("E433").Value = 84
$ws.Range("G433").Value = 0
$ws.Range("H433").Value = 0
$ws.Range("I433").Value = "No"
$ws.Range("J433").Value = 3
$ws.Range("K433").Value = "Neutral"
$ws.Range("L433").Value = 7
$ws.Range("M433").Value = 0.63
$ws.Range("N433").Value = 68
$ws.Range("O433").Value = 30.03
$ws.Range("P433").Value = 4
$ws.Range("Q433").Value = 0.31
$ws.Range("R433").Value = 9.9
$ws.Range("S433").Value = 57
$ws.Range("T433").Value = 0
$ws.Range("A434").Value = 45848
$ws.Range("A434").NumberFormat = "m/d/yy"
$ws.Range("B434").Value = "Nonflowering"
$ws.Range("C434").Value = "Large"
$ws.Range("D434").Value = 69
$ws.Range("E434").Value = 84
$ws.Range("G434").Value = 0
$ws.Range("H434").Value = 0.2
$ws.Range("I434").Value = "No"
$ws.Range("J434").Value = 4
$ws.Range("K434").Value = "Dark"
$ws.Range("L434").Value = 7
$ws.Range("M434").Value = 0.63
$ws.Range("N434").Value = 68
$ws.Range("O434").Value = 30.03
$ws.Range("P434").Value = 4
$ws.Range("Q434").Value = 0.31
$ws.Range("R434").Value = 9.9
$ws.Range("S434").Value = 57
$ws.Range("T434").Value = 0
$ws.Range("A435").Value = 45848
$ws.Range("A435").NumberFormat = "m/d/yy"
$ws.Range("B435").Value = "Tree"
$ws.Range("C435").Value = "Medium"
$ws.Range("D435").Value = 69
$ws.Range("E435").Value = 84
$ws.Range("G435").Value = 0
$ws.Range("H435").Value = 0.2
$ws.Range("I435").Value = "No"
$ws.Range("J435").Value = 1
$ws.Range("K435").Value = "Bright"
$ws.Range("L435").Value = 7
$ws.Range("M435").Value = 0.63
$ws.Range("N435").Value = 68
$ws.Range("O435").Value = 30.03
$ws.Range("P435").Value = 4
$ws.Range("Q435").Value = 0.31
$ws.Range("R435").Value = 9.9
$ws.Range("S435").Value = 57
$ws.Range("T435").Value = 0

# Temp_Diff column is a shared formula (=ABS(Low-High)); extend it over the new rows
$ws.Range("F429:F435").Formula = "=ABS(D429-E429)"

# Restore selection to reflect the new extent of the data (post-edit view state)
$ws.Range("I436").Select()
